$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a known-empty, default-styled cell as a format donor so that
# forcing numeric-looking Price strings to Text does not leave a
# lingering per-cell NumberFormat override behind.
$ws.Range("Z1").Copy() | Out-Null

$ws.Range("D2").Value = "69.168.73"
$ws.Range("E2").Value = "  +1.01%  "
$ws.Range("D3").Value = "3.775.12"
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "632.27"
$ws.Range("D5").PasteSpecial(-4122) | Out-Null
$ws.Range("E5").Value = "  +3.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.31"
$ws.Range("D6").PasteSpecial(-4122) | Out-Null
$ws.Range("E6").Value = "  +1.54%  "
$ws.Range("D7").Value = "3.773.22"
$ws.Range("E7").Value = "  -0.76%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("E9").Value = "  +0.83%  "
$ws.Range("E10").Value = "  -0.55%  "
$ws.Range("E11").Value = "  +2.27%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.79"
$ws.Range("D12").PasteSpecial(-4122) | Out-Null
$ws.Range("E12").Value = "  -2.51%  "
$ws.Range("E13").Value = "  -1.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.89"
$ws.Range("D14").PasteSpecial(-4122) | Out-Null
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("D15").Value = "4.409.08"
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("D16").Value = "3.777.06"
$ws.Range("E16").Value = "  -1.34%  "
$ws.Range("D17").Value = "69.214.52"
$ws.Range("E17").Value = "  +1.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.65"
$ws.Range("D18").PasteSpecial(-4122) | Out-Null
$ws.Range("E18").Value = "  -2.42%  "
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "463.75"
$ws.Range("D21").PasteSpecial(-4122) | Out-Null
$ws.Range("E21").Value = "  +0.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.55"
$ws.Range("D22").PasteSpecial(-4122) | Out-Null
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("E23").Value = "  +1.45%  "
$ws.Range("E24").Value = "  -0.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.79"
$ws.Range("D25").PasteSpecial(-4122) | Out-Null
$ws.Range("E25").Value = "  -0.83%  "
$ws.Range("E26").Value = "  +0.64%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.14"
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("E27").Value = "  +1.70%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.10"
$ws.Range("D28").PasteSpecial(-4122) | Out-Null
$ws.Range("E28").Value = "  +1.06%  "
$ws.Range("D30").Value = "3.925.45"
$ws.Range("E30").Value = "  -0.52%  "
$ws.Range("E32").Value = "  +2.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.11"
$ws.Range("D33").PasteSpecial(-4122) | Out-Null
$ws.Range("E33").Value = "  -1.89%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.52"
$ws.Range("D34").PasteSpecial(-4122) | Out-Null
$ws.Range("E34").Value = "  -1.83%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.169"
$ws.Range("D35").PasteSpecial(-4122) | Out-Null
$ws.Range("E35").Value = "  +14.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("D36").PasteSpecial(-4122) | Out-Null
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("D37").Value = "3.728.34"
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("E38").Value = "  -0.91%  "
$ws.Range("E39").Value = "  +0.38%  "
$ws.Range("E40").Value = "  +3.99%  "
$ws.Range("E41").Value = "  -1.42%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.965"
$ws.Range("D42").PasteSpecial(-4122) | Out-Null
$ws.Range("E42").Value = "  -1.65%  "
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "158.17"
$ws.Range("D45").PasteSpecial(-4122) | Out-Null
$ws.Range("E45").Value = "  +3.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.96"
$ws.Range("D46").PasteSpecial(-4122) | Out-Null
$ws.Range("E46").Value = "  +5.27%  "
$ws.Range("E47").Value = "  +1.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "43.20"
$ws.Range("D48").PasteSpecial(-4122) | Out-Null
$ws.Range("E48").Value = "  +0.32%  "
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "46.70"
$ws.Range("D50").PasteSpecial(-4122) | Out-Null
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.38"
$ws.Range("D51").PasteSpecial(-4122) | Out-Null
$ws.Range("E51").Value = "  -0.01%  "

$excel.CutCopyMode = 0
